# Refactor the "label / value" two-column layout used on every slide so
# that it is more failsafe: the label column is narrower (and shifted
# slightly left), the value column starts a bit further left, and every
# row is taller (24pt instead of 20pt) with rows re-flowed accordingly.

$p = $ppt.ActivePresentation

# Column geometry (points). Label = left column (e.g. "name:"), Value =
# right column (e.g. "Victoria Hanover").
$labelLeft  = 235
$labelWidth = 115
$valueLeft  = 350
$valueWidth = 350
$rowHeight  = 24
$topStart   = 20

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)

        # Rows alternate label (odd shape index) / value (even shape index).
        $row = [Math]::Floor(($i - 1) / 2)
        $top = $topStart + ($row * $rowHeight)

        if (($i % 2) -eq 1) {
            $sh.Left   = $labelLeft
            $sh.Width  = $labelWidth
        } else {
            $sh.Left   = $valueLeft
            $sh.Width  = $valueWidth
        }

        $sh.Top    = $top
        $sh.Height = $rowHeight
    }
}
